$d = $word.ActiveDocument
$d.Content.Find.Execute("more than 250 tool supporting it", $true, $false, $false, $false, $false,
                         $true, 1, $false, "more than 250 tools supporting it", 2)
